# Auto-generated edit script: apply value updates from the source diff
# to the corresponding worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 109.22222
$ws.Range("I33").Value = 103.69231
$ws.Range("J33").Value = 123.6
$ws.Range("K33").Value = 103.69231
$ws.Range("L33").Value = 123.6
$ws.Range("M33").Value = 125.30769
$ws.Range("N33").Value = -581.6
$ws.Range("H105").Value = 17500
$ws.Range("J105").Value = 17500
$ws.Range("L105").Value = 17500
$ws.Range("N105").Value = -24488
$ws.Range("H107").Value = 181.2
$ws.Range("I107").Value = 186
$ws.Range("J107").Value = 174
$ws.Range("K107").Value = 186
$ws.Range("L107").Value = 174
$ws.Range("M107").Value = 1734
$ws.Range("N107").Value = -4014
$ws.Range("H127").Value = 1460.6666
$ws.Range("I127").Value = 768.25
$ws.Range("K127").Value = 2304.75
$ws.Range("M127").Value = 2655.25
$ws.Range("H141").Value = 2670.3333
$ws.Range("I141").Value = 2850.0908
$ws.Range("J141").Value = 693
$ws.Range("K141").Value = 8550.2724
$ws.Range("L141").Value = 2079
$ws.Range("M141").Value = -3370.2724
$ws.Range("N141").Value = -12439

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1317.2222
$ws.Range("I132").Value = 1198.5217
$ws.Range("K132").Value = 3595.5651
$ws.Range("M132").Value = -1065.5651

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 20279.334
$ws.Range("J76").Value = 20279.334
$ws.Range("L76").Value = 20279.334
$ws.Range("N76").Value = -20909.334
$ws.Range("H79").Value = 20279.334
$ws.Range("J79").Value = 20279.334
$ws.Range("L79").Value = 20279.334
$ws.Range("N79").Value = -22463.334
$ws.Range("H94").Value = 2470.125
$ws.Range("I94").Value = 2470.125
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2470.125
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2019.125
$ws.Range("H134").Value = 2476.6
$ws.Range("I134").Value = 2476.6
$ws.Range("K134").Value = 7429.799999999999
$ws.Range("M134").Value = -4894.799999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 342
$ws.Range("J2").Value = 252.33333
$ws.Range("L2").Value = 252.33333
$ws.Range("N2").Value = -478.33333
$ws.Range("H3").Value = 3668333.2
$ws.Range("J3").Value = 502500
$ws.Range("L3").Value = 502500
$ws.Range("N3").Value = -502726
$ws.Range("H5").Value = 152.36363
$ws.Range("J5").Value = 97.666664
$ws.Range("L5").Value = 97.666664
$ws.Range("N5").Value = -321.666664
$ws.Range("H7").Value = 187.55556
$ws.Range("I7").Value = 117.6
$ws.Range("J7").Value = 275
$ws.Range("K7").Value = 117.6
$ws.Range("L7").Value = 275
$ws.Range("M7").Value = -4.599999999999994
$ws.Range("N7").Value = -501
$ws.Range("H8").Value = 3261.8
$ws.Range("I8").Value = 2700
$ws.Range("J8").Value = 3636.3333
$ws.Range("K8").Value = 2700
$ws.Range("L8").Value = 3636.3333
$ws.Range("M8").Value = -2560
$ws.Range("N8").Value = -3916.3333
$ws.Range("H12").Value = 4117.625
$ws.Range("I12").Value = 612
$ws.Range("J12").Value = 7623.25
$ws.Range("K12").Value = 612
$ws.Range("L12").Value = 7623.25
$ws.Range("M12").Value = -442
$ws.Range("N12").Value = -7963.25
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H15").Value = 73996.5
$ws.Range("J15").Value = 73996.5
$ws.Range("L15").Value = 73996.5
$ws.Range("N15").Value = -74336.5
$ws.Range("H16").Value = 1548.7
$ws.Range("I16").Value = 1639.8572
$ws.Range("K16").Value = 1639.8572
$ws.Range("M16").Value = -1352.8572
$ws.Range("H31").Value = 2225.2104
$ws.Range("I31").Value = 1668.0385
$ws.Range("K31").Value = 1668.0385
$ws.Range("M31").Value = -1373.0385
$ws.Range("H34").Value = 2225.2104
$ws.Range("I34").Value = 1668.0385
$ws.Range("K34").Value = 1668.0385
$ws.Range("M34").Value = -1466.0385
$ws.Range("H58").Value = 5957
$ws.Range("I58").Value = 5350.4614
$ws.Range("K58").Value = 5350.4614
$ws.Range("M58").Value = -5147.4614
$ws.Range("H107").Value = 1183.1111
$ws.Range("I107").Value = 324.5
$ws.Range("J107").Value = 1428.4286
$ws.Range("K107").Value = 324.5
$ws.Range("L107").Value = 1428.4286
$ws.Range("M107").Value = 1595.5
$ws.Range("N107").Value = -5268.4286
$ws.Range("H113").Value = 1548.7
$ws.Range("I113").Value = 1639.8572
$ws.Range("K113").Value = 1639.8572
$ws.Range("M113").Value = 530.1428000000001
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139
$ws.Range("H136").Value = 5957
$ws.Range("I136").Value = 5350.4614
$ws.Range("K136").Value = 16051.3842
$ws.Range("M136").Value = -13501.3842

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 689.4
$ws.Range("I60").Value = 386.75
$ws.Range("J60").Value = 1900
$ws.Range("K60").Value = 1160.25
$ws.Range("L60").Value = 5700
$ws.Range("M60").Value = -909.25
$ws.Range("N60").Value = -6202
$ws.Range("H63").Value = 1178.6
$ws.Range("J63").Value = 1499.5
$ws.Range("L63").Value = 4498.5
$ws.Range("N63").Value = -5996.5
$ws.Range("H66").Value = 1178.6
$ws.Range("J66").Value = 1499.5
$ws.Range("L66").Value = 13495.5
$ws.Range("N66").Value = -20983.5
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H87").Value = 400
$ws.Range("I87").Value = 400
$ws.Range("K87").Value = 1200
$ws.Range("M87").Value = 48
$ws.Range("H90").Value = 400
$ws.Range("I90").Value = 400
$ws.Range("K90").Value = 3600
$ws.Range("M90").Value = 2640
$ws.Range("H103").Value = 139.57143
$ws.Range("I103").Value = 154.5
$ws.Range("K103").Value = 463.5
$ws.Range("M103").Value = 415.5
$ws.Range("H132").Value = 629.3333
$ws.Range("I132").Value = 629.3333
$ws.Range("K132").Value = 5663.9997
$ws.Range("M132").Value = -3133.9997
$ws.Range("H139").Value = 3377.348
$ws.Range("I139").Value = 853.2222
$ws.Range("K139").Value = 2559.6666
$ws.Range("M139").Value = 2580.3334
$ws.Range("H140").Value = 2425
$ws.Range("I140").Value = 2425
$ws.Range("K140").Value = 7275
$ws.Range("M140").Value = -2095

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1366
$ws.Range("I132").Value = 1366
$ws.Range("K132").Value = 4098
$ws.Range("M132").Value = -1568

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 992.875
$ws.Range("I22").Value = 992.875
$ws.Range("K22").Value = 992.875
$ws.Range("M22").Value = -697.875
$ws.Range("H27").Value = 992.875
$ws.Range("I27").Value = 992.875
$ws.Range("K27").Value = 992.875
$ws.Range("M27").Value = -885.875
$ws.Range("H46").Value = 3659.3333
$ws.Range("I46").Value = 2481.6667
$ws.Range("J46").Value = 4444.4443
$ws.Range("K46").Value = 2481.6667
$ws.Range("L46").Value = 4444.4443
$ws.Range("M46").Value = -2293.6667
$ws.Range("N46").Value = -4820.4443

Write-Output "Applied all Marilith_Profits cell updates"
